$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J17").Value = 4792.4
$ws.Range("N17").Value = -14713.2
$ws.Range("L17").Value = 14377.2
$ws.Range("H17").Value = 4810
$ws.Range("I39").Value = 198.83333
$ws.Range("H39").Value = 256.14285
$ws.Range("M39").Value = -300.49999
$ws.Range("K39").Value = 596.49999
$ws.Range("N40").Value = -4341.5833
$ws.Range("J40").Value = 3991.5833
$ws.Range("H40").Value = 3987.92
$ws.Range("L40").Value = 3991.5833
$ws.Range("M43").Value = -2958.6667
$ws.Range("N43").Value = -4797.3
$ws.Range("I43").Value = 3027.6667
$ws.Range("K43").Value = 3027.6667
$ws.Range("H43").Value = 4047.4375
$ws.Range("J43").Value = 4659.3
$ws.Range("L43").Value = 4659.3
$ws.Range("H58").Value = 2465.1428
$ws.Range("J58").Value = 2177.6667
$ws.Range("L58").Value = 6533.000100000001
$ws.Range("N58").Value = -6833.000100000001
$ws.Range("M62").Value = -4625.2
$ws.Range("H62").Value = 5951.7646
$ws.Range("K62").Value = 5249.2
$ws.Range("L62").Value = 6955.4287
$ws.Range("J62").Value = 6955.4287
$ws.Range("N62").Value = -8203.4287
$ws.Range("I62").Value = 5249.2
$ws.Range("H64").Value = 7219.9
$ws.Range("I64").Value = 2599
$ws.Range("M64").Value = -2351
$ws.Range("K64").Value = 2599
$ws.Range("K65").Value = 26246
$ws.Range("M65").Value = -23126
$ws.Range("H65").Value = 5951.7646
$ws.Range("J65").Value = 6955.4287
$ws.Range("I65").Value = 5249.2
$ws.Range("N65").Value = -41017.14350000001
$ws.Range("L65").Value = 34777.14350000001
$ws.Range("H67").Value = 7219.9
$ws.Range("M67").Value = -1741
$ws.Range("K67").Value = 2599
$ws.Range("I67").Value = 2599
$ws.Range("I76").Value = 4645.6665
$ws.Range("H76").Value = 5017.3
$ws.Range("K76").Value = 4645.6665
$ws.Range("M76").Value = -4330.6665
$ws.Range("M79").Value = -3553.6665
$ws.Range("H79").Value = 5017.3
$ws.Range("K79").Value = 4645.6665
$ws.Range("I79").Value = 4645.6665
$ws.Range("K80").Value = 500004630
$ws.Range("M80").Value = -500003632
$ws.Range("H80").Value = 100019464
$ws.Range("I80").Value = 166668210
$ws.Range("M83").Value = -1500008898
$ws.Range("K83").Value = 1500013890
$ws.Range("H83").Value = 100019464
$ws.Range("I83").Value = 166668210
$ws.Range("I86").Value = 333336130
$ws.Range("H86").Value = 200001870
$ws.Range("K86").Value = 333336130
$ws.Range("M86").Value = -333335007
$ws.Range("H88").Value = 2193.75
$ws.Range("J88").Value = 1487.5
$ws.Range("L88").Value = 1487.5
$ws.Range("N88").Value = -2299.5
$ws.Range("I89").Value = 333336130
$ws.Range("H89").Value = 200001870
$ws.Range("K89").Value = 1666680650
$ws.Range("M89").Value = -1666675034
$ws.Range("H91").Value = 2193.75
$ws.Range("J91").Value = 1487.5
$ws.Range("L91").Value = 1487.5
$ws.Range("N91").Value = -4295.5
$ws.Range("N95").ClearContents()
$ws.Range("J95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M98").Value = 42.44450000000006
$ws.Range("H98").Value = 1509.7
$ws.Range("K98").Value = 1455.5555
$ws.Range("I98").Value = 1455.5555
$ws.Range("L111").Value = 11134.2
$ws.Range("H111").Value = 3652.4443
$ws.Range("J111").Value = 3711.4
$ws.Range("N111").Value = -17268.2
$ws.Range("H113").Value = 22414
$ws.Range("J113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("L113").Value = 0
$ws.Range("H117").Value = 79100.836
$ws.Range("L117").Value = 79100.836
$ws.Range("N117").Value = -88278.836
$ws.Range("J117").Value = 79100.836
$ws.Range("H121").Value = 2380.55
$ws.Range("J121").Value = 2380.55
$ws.Range("L121").Value = 7141.650000000001
$ws.Range("N121").Value = -10635.65
$ws.Range("I122").Value = 1455.5555
$ws.Range("K122").Value = 4366.666499999999
$ws.Range("M122").Value = -1916.666499999999
$ws.Range("H122").Value = 1509.7
$ws.Range("L129").Value = 6599.571599999999
$ws.Range("N129").Value = -16599.5716
$ws.Range("H129").Value = 1763.4546
$ws.Range("J129").Value = 2199.8572
$ws.Range("N132").Value = -18719.8568
$ws.Range("K132").Value = 5902.559999999999
$ws.Range("M132").Value = -3372.559999999999
$ws.Range("J132").Value = 4553.2856
$ws.Range("I132").Value = 1967.52
$ws.Range("H132").Value = 2533.1562
$ws.Range("L132").Value = 13659.8568
$ws.Range("I135").Value = 920.6111
$ws.Range("L135").Value = 54672.75
$ws.Range("M135").Value = -5750.499899999999
$ws.Range("J135").Value = 6074.75
$ws.Range("K135").Value = 8285.499899999999
$ws.Range("N135").Value = -59742.75
$ws.Range("H135").Value = 1857.7273
$ws.Range("M138").Value = 710.5
$ws.Range("J138").Value = 2733.8936
$ws.Range("K138").Value = 4429.5
$ws.Range("N138").Value = -18481.6808
$ws.Range("I138").Value = 1476.5
$ws.Range("H138").Value = 2658.45
$ws.Range("L138").Value = 8201.6808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J23").Value = 20771.666
$ws.Range("N23").Value = -21289.666
$ws.Range("H23").Value = 20771.666
$ws.Range("L23").Value = 20771.666
$ws.Range("I26").Value = 7764.5713
$ws.Range("H26").Value = 7764.5713
$ws.Range("M26").Value = -7434.5713
$ws.Range("K26").Value = 7764.5713
$ws.Range("K32").Value = 26582086
$ws.Range("H32").Value = 24419464
$ws.Range("I32").Value = 26582086
$ws.Range("M32").Value = -26581799
$ws.Range("H35").Value = 10478.667
$ws.Range("K35").Value = 1718
$ws.Range("M35").Value = -1312
$ws.Range("I35").Value = 1718
$ws.Range("M45").Value = -4451.222
$ws.Range("I45").Value = 4828.222
$ws.Range("K45").Value = 4828.222
$ws.Range("H45").Value = 4585.8184
$ws.Range("J61").Value = 4068.625
$ws.Range("H61").Value = 3642.4119
$ws.Range("L61").Value = 4068.625
$ws.Range("N61").Value = -4492.625
$ws.Range("I61").Value = 3263.5557
$ws.Range("M61").Value = -3051.5557
$ws.Range("K61").Value = 3263.5557
$ws.Range("K74").Value = 3681
$ws.Range("I74").Value = 3681
$ws.Range("H74").Value = 3393.2104
$ws.Range("M74").Value = -2807
$ws.Range("H77").Value = 3393.2104
$ws.Range("I77").Value = 3681
$ws.Range("K77").Value = 18405
$ws.Range("M77").Value = -14037
$ws.Range("H88").Value = 1642.8334
$ws.Range("K88").Value = 1450
$ws.Range("M88").Value = -1044
$ws.Range("I88").Value = 1450
$ws.Range("J88").Value = 1739.25
$ws.Range("L88").Value = 1739.25
$ws.Range("N88").Value = -2551.25
$ws.Range("H91").Value = 1642.8334
$ws.Range("J91").Value = 1739.25
$ws.Range("L91").Value = 1739.25
$ws.Range("K91").Value = 1450
$ws.Range("M91").Value = -46
$ws.Range("N91").Value = -4547.25
$ws.Range("I91").Value = 1450
$ws.Range("K102").Value = 1612.7333
$ws.Range("H102").Value = 2062.6843
$ws.Range("J102").Value = 3750
$ws.Range("N102").Value = -6994
$ws.Range("M102").Value = 9.266699999999901
$ws.Range("L102").Value = 3750
$ws.Range("I102").Value = 1612.7333
$ws.Range("I122").Value = 3933.3333
$ws.Range("L122").Value = 20000.6661
$ws.Range("N122").Value = -24900.6661
$ws.Range("K122").Value = 11799.9999
$ws.Range("M122").Value = -9349.999899999999
$ws.Range("H122").Value = 5300.1113
$ws.Range("J122").Value = 6666.8887
$ws.Range("N130").ClearContents()
$ws.Range("L130").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("N132").Value = -20288
$ws.Range("K132").Value = 9236.832900000001
$ws.Range("M132").Value = -6706.832900000001
$ws.Range("J132").Value = 5076
$ws.Range("I132").Value = 3078.9443
$ws.Range("H132").Value = 3442.0454
$ws.Range("L132").Value = 15228
$ws.Range("L136").Value = 12205.875
$ws.Range("N136").Value = -17305.875
$ws.Range("J136").Value = 4068.625
$ws.Range("M136").Value = -7240.667099999999
$ws.Range("H136").Value = 3642.4119
$ws.Range("I136").Value = 3263.5557
$ws.Range("K136").Value = 9790.667099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 1351.2222
$ws.Range("H20").Value = 1738.1428
$ws.Range("K20").Value = 1351.2222
$ws.Range("N20").Value = -2928.6
$ws.Range("L20").Value = 2434.6
$ws.Range("M20").Value = -1104.2222
$ws.Range("J20").Value = 2434.6
$ws.Range("N63").Value = -85391.375
$ws.Range("H63").Value = 84019.375
$ws.Range("J63").Value = 84019.375
$ws.Range("L63").Value = 84019.375
$ws.Range("N66").Value = -258922.125
$ws.Range("L66").Value = 252058.125
$ws.Range("J66").Value = 84019.375
$ws.Range("H66").Value = 84019.375
$ws.Range("I86").Value = 5253
$ws.Range("H86").Value = 3835.3333
$ws.Range("K86").Value = 5253
$ws.Range("M86").Value = -4130
$ws.Range("I89").Value = 5253
$ws.Range("H89").Value = 3835.3333
$ws.Range("K89").Value = 26265
$ws.Range("M89").Value = -20649
$ws.Range("J94").Value = 987.8182
$ws.Range("H94").Value = 1086.7
$ws.Range("L94").Value = 987.8182
$ws.Range("M94").Value = -756.5554999999999
$ws.Range("N94").Value = -1889.8182
$ws.Range("K94").Value = 1207.5555
$ws.Range("I94").Value = 1207.5555
$ws.Range("I105").Value = 1132.9231
$ws.Range("H105").Value = 1981.8096
$ws.Range("K105").Value = 1132.9231
$ws.Range("M105").Value = 614.0769
$ws.Range("H107").Value = 2750.6667
$ws.Range("M107").Value = 916.4
$ws.Range("K107").Value = 1003.6
$ws.Range("I107").Value = 1003.6
$ws.Range("K134").Value = 8247462
$ws.Range("M134").Value = -8244927
$ws.Range("H134").Value = 2465324.2
$ws.Range("I134").Value = 2749154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -1550.5
$ws.Range("H16").Value = 1837.5
$ws.Range("I16").Value = 1837.5
$ws.Range("K16").Value = 1837.5
$ws.Range("H69").Value = 58990.25
$ws.Range("M69").Value = -19433
$ws.Range("K69").Value = 20182
$ws.Range("I69").Value = 20182
$ws.Range("H72").Value = 58990.25
$ws.Range("M72").Value = -56802
$ws.Range("I72").Value = 20182
$ws.Range("K72").Value = 60546
$ws.Range("L81").Value = 41664
$ws.Range("H81").Value = 44442.332
$ws.Range("J81").Value = 41664
$ws.Range("N81").Value = -43660
$ws.Range("N84").Value = -134976
$ws.Range("J84").Value = 41664
$ws.Range("H84").Value = 44442.332
$ws.Range("L84").Value = 124992
$ws.Range("N87").Value = -99341
$ws.Range("L87").Value = 96969
$ws.Range("H87").Value = 96969
$ws.Range("J87").Value = 96969
$ws.Range("H88").Value = 29383.715
$ws.Range("J88").Value = 29383.715
$ws.Range("L88").Value = 29383.715
$ws.Range("N88").Value = -30195.715
$ws.Range("N90").Value = -302763
$ws.Range("J90").Value = 96969
$ws.Range("H90").Value = 96969
$ws.Range("L90").Value = 290907
$ws.Range("H91").Value = 29383.715
$ws.Range("J91").Value = 29383.715
$ws.Range("L91").Value = 29383.715
$ws.Range("N91").Value = -32191.715
$ws.Range("I105").Value = 3468.5
$ws.Range("H105").Value = 3481.2856
$ws.Range("K105").Value = 3468.5
$ws.Range("M105").Value = -1721.5
$ws.Range("K113").Value = 1837.5
$ws.Range("M113").Value = 332.5
$ws.Range("H113").Value = 1837.5
$ws.Range("I113").Value = 1837.5
$ws.Range("I122").Value = 4749
$ws.Range("K122").Value = 14247
$ws.Range("M122").Value = -11797
$ws.Range("H122").Value = 5298.75
$ws.Range("N132").Value = -20401.571
$ws.Range("K132").Value = 12060.4284
$ws.Range("M132").Value = -9530.428400000001
$ws.Range("J132").Value = 5113.857
$ws.Range("I132").Value = 4020.1428
$ws.Range("H132").Value = 4293.5713
$ws.Range("L132").Value = 15341.571
$ws.Range("K134").Value = 6193.928400000001
$ws.Range("M134").Value = -3658.928400000001
$ws.Range("H134").Value = 2482.889
$ws.Range("I134").Value = 2064.6428
$ws.Range("N141").Value = -500924
$ws.Range("J141").Value = 490564
$ws.Range("L141").Value = 490564
$ws.Range("H141").Value = 490564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I14").Value = 14295.125
$ws.Range("M14").Value = -42712.375
$ws.Range("K14").Value = 42885.375
$ws.Range("H14").Value = 14295.125
$ws.Range("J23").Value = 166.66667
$ws.Range("I23").Value = 170
$ws.Range("M23").Value = -275
$ws.Range("K23").Value = 510
$ws.Range("N23").Value = -970.00001
$ws.Range("H23").Value = 168
$ws.Range("L23").Value = 500.00001
$ws.Range("K33").Value = 3639
$ws.Range("I33").Value = 606.5
$ws.Range("H33").Value = 979.2353000000001
$ws.Range("M33").Value = -3356
$ws.Range("I86").Value = 317.6
$ws.Range("H86").Value = 385.375
$ws.Range("K86").Value = 952.8000000000001
$ws.Range("M86").Value = 233.1999999999999
$ws.Range("I89").Value = 317.6
$ws.Range("H89").Value = 385.375
$ws.Range("K89").Value = 2858.4
$ws.Range("M89").Value = 3069.6
$ws.Range("H107").Value = 793.3889
$ws.Range("N107").Value = -5530.95
$ws.Range("L107").Value = 1690.95
$ws.Range("M107").Value = -1321.6875
$ws.Range("K107").Value = 3241.6875
$ws.Range("I107").Value = 1080.5625
$ws.Range("J107").Value = 563.65
$ws.Range("J118").Value = 4000
$ws.Range("H118").Value = 2276.3333
$ws.Range("I118").Value = 1414.5
$ws.Range("M118").Value = -3000.5
$ws.Range("N118").Value = -14486
$ws.Range("K118").Value = 4243.5
$ws.Range("L118").Value = 12000
$ws.Range("I122").Value = 0
$ws.Range("L122").Value = 130392
$ws.Range("N122").Value = -135292
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 14488
$ws.Range("J122").Value = 14488
$ws.Range("N123").ClearContents()
$ws.Range("M123").Value = -14017
$ws.Range("I123").Value = 5489
$ws.Range("H123").Value = 5489
$ws.Range("K123").Value = 16467
$ws.Range("L123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N132").Value = -17316.2
$ws.Range("K132").Value = 4920.75
$ws.Range("M132").Value = -2390.75
$ws.Range("J132").Value = 1361.8
$ws.Range("I132").Value = 546.75
$ws.Range("H132").Value = 999.55554
$ws.Range("L132").Value = 12256.2
$ws.Range("L140").Value = 15000
$ws.Range("N140").Value = -25360
$ws.Range("M140").Value = -365.3335000000006
$ws.Range("K140").Value = 5545.333500000001
$ws.Range("H140").Value = 2636.3333
$ws.Range("J140").Value = 5000
$ws.Range("I140").Value = 1848.4445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L52").Value = 45070
$ws.Range("J52").Value = 45070
$ws.Range("N52").Value = -45588
$ws.Range("H52").Value = 45455
$ws.Range("K80").Value = 3340
$ws.Range("M80").Value = -2342
$ws.Range("H80").Value = 3485.7144
$ws.Range("I80").Value = 3340
$ws.Range("M83").Value = -11708
$ws.Range("K83").Value = 16700
$ws.Range("H83").Value = 3485.7144
$ws.Range("I83").Value = 3340
$ws.Range("N95").ClearContents()
$ws.Range("J95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N96").Value = -80491.5
$ws.Range("J96").Value = 74999.5
$ws.Range("H96").Value = 74999.5
$ws.Range("L96").Value = 74999.5
$ws.Range("I97").Value = 1769.4
$ws.Range("H97").Value = 2365.6
$ws.Range("L97").Value = 2961.8
$ws.Range("K97").Value = 1769.4
$ws.Range("N97").Value = -3953.8
$ws.Range("J97").Value = 2961.8
$ws.Range("M97").Value = -1273.4
$ws.Range("H117").Value = 77489.5
$ws.Range("L117").Value = 77489.5
$ws.Range("N117").Value = -84373.5
$ws.Range("J117").Value = 77489.5
$ws.Range("I122").Value = 1594
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("K122").Value = 4782
$ws.Range("M122").Value = -2332
$ws.Range("H122").Value = 1797
$ws.Range("J122").Value = 2000
$ws.Range("J126").Value = 3911.6667
$ws.Range("I126").Value = 3910
$ws.Range("M126").Value = -9260
$ws.Range("K126").Value = 11730
$ws.Range("N126").Value = -16675.0001
$ws.Range("L126").Value = 11735.0001
$ws.Range("H126").Value = 3911
$ws.Range("K132").Value = 11760.9999
$ws.Range("M132").Value = -9230.999899999999
$ws.Range("I132").Value = 3920.3333
$ws.Range("H132").Value = 3926.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -13388
$ws.Range("H7").Value = 11318
$ws.Range("I7").Value = 13500
$ws.Range("K7").Value = 13500
$ws.Range("I22").Value = 1791.1
$ws.Range("K22").Value = 1791.1
$ws.Range("H22").Value = 3038.5
$ws.Range("M22").Value = -1496.1
$ws.Range("M27").Value = -1684.1
$ws.Range("I27").Value = 1791.1
$ws.Range("K27").Value = 1791.1
$ws.Range("H27").Value = 3038.5
$ws.Range("M31").Value = -1272.8334
$ws.Range("K31").Value = 1520.8334
$ws.Range("I31").Value = 1520.8334
$ws.Range("H31").Value = 2603.8462
$ws.Range("N40").Value = -5007
$ws.Range("J40").Value = 4735
$ws.Range("H40").Value = 4338
$ws.Range("L40").Value = 4735
$ws.Range("K46").Value = 2938.2
$ws.Range("H46").Value = 8108.76
$ws.Range("M46").Value = -2750.2
$ws.Range("I46").Value = 2938.2
$ws.Range("N82").Value = -2810
$ws.Range("I82").Value = 1532.875
$ws.Range("M82").Value = -1171.875
$ws.Range("K82").Value = 1532.875
$ws.Range("J82").Value = 2088
$ws.Range("H82").Value = 1746.3846
$ws.Range("L82").Value = 2088
$ws.Range("I85").Value = 1532.875
$ws.Range("K85").Value = 1532.875
$ws.Range("L85").Value = 2088
$ws.Range("M85").Value = -284.875
$ws.Range("N85").Value = -4584
$ws.Range("H85").Value = 1746.3846
$ws.Range("J85").Value = 2088
$ws.Range("H100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("I100").Value = 2000
$ws.Range("N115").Value = -57339
$ws.Range("H115").Value = 54989
$ws.Range("J115").Value = 54989
$ws.Range("L115").Value = 54989
$ws.Range("I126").Value = 13500
$ws.Range("M126").Value = -38030
$ws.Range("K126").Value = 40500
$ws.Range("H126").Value = 11318
$ws.Range("N132").Value = -20516.231
$ws.Range("K132").Value = 12491.4552
$ws.Range("M132").Value = -9961.4552
$ws.Range("J132").Value = 5152.077
$ws.Range("I132").Value = 4163.8184
$ws.Range("H132").Value = 4530.8857
$ws.Range("L132").Value = 15456.231
$ws.Range("M136").Value = -12394.6362
$ws.Range("H136").Value = 7610.2354
$ws.Range("I136").Value = 4981.5454
$ws.Range("K136").Value = 14944.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value = 1622
$ws.Range("M81").Value = -2183
$ws.Range("H81").Value = 2715.8333
$ws.Range("K81").Value = 3244
$ws.Range("I84").Value = 1622
$ws.Range("K84").Value = 16220
$ws.Range("M84").Value = -10916
$ws.Range("H84").Value = 2715.8333
$ws.Range("K113").Value = 1155
$ws.Range("M113").Value = 1015
$ws.Range("H113").Value = 640
$ws.Range("J113").Value = 1150
$ws.Range("I113").Value = 385
$ws.Range("N113").Value = -7790
$ws.Range("L113").Value = 3450
$ws.Range("N120").Value = -90623.42999999999
$ws.Range("H120").Value = 80947.42999999999
$ws.Range("L120").Value = 80947.42999999999
$ws.Range("J120").Value = 80947.42999999999
$ws.Range("I122").Value = 13598.75
$ws.Range("L122").Value = 40495.5
$ws.Range("N122").Value = -45395.5
$ws.Range("K122").Value = 40796.25
$ws.Range("M122").Value = -38346.25
$ws.Range("H122").Value = 13565.333
$ws.Range("J122").Value = 13498.5
$ws.Range("H125").Value = 131999
$ws.Range("J125").Value = 131999
$ws.Range("L125").Value = 131999
$ws.Range("N125").Value = -141839
$ws.Range("J126").Value = 3000
$ws.Range("I126").Value = 4500
$ws.Range("M126").Value = -11030
$ws.Range("K126").Value = 13500
$ws.Range("N126").Value = -13940
$ws.Range("L126").Value = 9000
$ws.Range("H126").Value = 4200
$ws.Range("K132").Value = 5269.5
$ws.Range("M132").Value = -2739.5
$ws.Range("I132").Value = 1756.5
$ws.Range("H132").Value = 2003.9286
$ws.Range("N141").Value = -90360
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("H141").Value = 80000
